$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 704.8333
$ws.Range("I2").Value = 557.5
$ws.Range("K2").Value = 557.5
$ws.Range("M2").Value = -444.5
$ws.Range("H12").Value = 1099.3334
$ws.Range("I12").Value = 556.2857
$ws.Range("K12").Value = 556.2857
$ws.Range("M12").Value = -386.2857
$ws.Range("H33").Value = 280
$ws.Range("I33").Value = 235.6
$ws.Range("K33").Value = 235.6
$ws.Range("M33").Value = -6.599999999999994
$ws.Range("H98").Value = 3820.8125
$ws.Range("I98").Value = 3555.5334
$ws.Range("K98").Value = 3555.5334
$ws.Range("M98").Value = -2057.5334
$ws.Range("H122").Value = 3820.8125
$ws.Range("I122").Value = 3555.5334
$ws.Range("K122").Value = 10666.6002
$ws.Range("M122").Value = -8216.600199999999
$ws.Range("H125").Value = 880.36365
$ws.Range("I125").Value = 882.6667
$ws.Range("J125").Value = 879.5
$ws.Range("K125").Value = 7944.0003
$ws.Range("L125").Value = 7915.5
$ws.Range("M125").Value = -5484.0003
$ws.Range("N125").Value = -12835.5
$ws.Range("H135").Value = 697
$ws.Range("I135").Value = 385.22223
$ws.Range("K135").Value = 3467.00007
$ws.Range("M135").Value = -932.0000700000001
$ws.Range("H137").Value = 2781054
$ws.Range("I137").Value = 3573100.5
$ws.Range("K137").Value = 10719301.5
$ws.Range("M137").Value = -10716751.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 839.3333
$ws.Range("I12").Value = 839.3333
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 839.3333
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -666.3333
$ws.Range("N12").ClearContents()
$ws.Range("H15").Value = 499
$ws.Range("J15").Value = 499
$ws.Range("L15").Value = 499
$ws.Range("N15").Value = -1199
$ws.Range("H32").Value = 1883294.9
$ws.Range("I32").Value = 843845.9399999999
$ws.Range("J32").Value = 47619050
$ws.Range("K32").Value = 843845.9399999999
$ws.Range("L32").Value = 47619050
$ws.Range("M32").Value = -843558.9399999999
$ws.Range("N32").Value = -47619624
$ws.Range("H74").Value = 2616.4075
$ws.Range("I74").Value = 2424.889
$ws.Range("J74").Value = 2999.4443
$ws.Range("K74").Value = 2424.889
$ws.Range("L74").Value = 2999.4443
$ws.Range("M74").Value = -1550.889
$ws.Range("N74").Value = -4747.4443
$ws.Range("H77").Value = 2616.4075
$ws.Range("I77").Value = 2424.889
$ws.Range("J77").Value = 2999.4443
$ws.Range("K77").Value = 12124.445
$ws.Range("L77").Value = 14997.2215
$ws.Range("M77").Value = -7756.445
$ws.Range("N77").Value = -23733.2215
$ws.Range("H125").Value = 124463.75
$ws.Range("J125").Value = 124463.75
$ws.Range("L125").Value = 124463.75
$ws.Range("N125").Value = -134303.75
$ws.Range("H132").Value = 4458.1665
$ws.Range("I132").Value = 4687.25
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 14061.75
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -11531.75
$ws.Range("N132").Value = -17060

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2895.8635
$ws.Range("I134").Value = 2179.6667
$ws.Range("K134").Value = 6539.000100000001
$ws.Range("M134").Value = -4004.000100000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2276224.8
$ws.Range("I31").Value = 2771.361
$ws.Range("K31").Value = 2771.361
$ws.Range("M31").Value = -2476.361
$ws.Range("H34").Value = 2276224.8
$ws.Range("I34").Value = 2771.361
$ws.Range("K34").Value = 2771.361
$ws.Range("M34").Value = -2569.361
$ws.Range("H86").Value = 6773.8096
$ws.Range("I86").Value = 6802.6113
$ws.Range("K86").Value = 6802.6113
$ws.Range("M86").Value = -5679.6113
$ws.Range("H89").Value = 6773.8096
$ws.Range("I89").Value = 6802.6113
$ws.Range("K89").Value = 34013.0565
$ws.Range("M89").Value = -28397.0565
$ws.Range("H99").Value = 3139.7144
$ws.Range("I99").Value = 1993.3334
$ws.Range("K99").Value = 1993.3334
$ws.Range("M99").Value = -495.3334
$ws.Range("H107").Value = 2942167.8
$ws.Range("J107").Value = 1670.625
$ws.Range("L107").Value = 1670.625
$ws.Range("N107").Value = -5510.625
$ws.Range("H126").Value = 3139.7144
$ws.Range("I126").Value = 1993.3334
$ws.Range("K126").Value = 5980.0002
$ws.Range("M126").Value = -3510.0002
$ws.Range("H132").Value = 16672522
$ws.Range("I132").Value = 3926.25
$ws.Range("K132").Value = 11778.75
$ws.Range("M132").Value = -9248.75
$ws.Range("H134").Value = 2341.087
$ws.Range("I134").Value = 2092.6904
$ws.Range("J134").Value = 4949.25
$ws.Range("K134").Value = 6278.0712
$ws.Range("L134").Value = 14847.75
$ws.Range("M134").Value = -3743.0712
$ws.Range("N134").Value = -19917.75

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 3900
$ws.Range("I102").Value = 3900
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 11700
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -9266
$ws.Range("N102").ClearContents()
$ws.Range("H121").Value = 5060949
$ws.Range("I121").Value = 12513308
$ws.Range("J121").Value = 92709
$ws.Range("K121").Value = 37539924
$ws.Range("L121").Value = 278127
$ws.Range("M121").Value = -37538614
$ws.Range("N121").Value = -280747
$ws.Range("H129").Value = 57212
$ws.Range("I129").Value = 765.6
$ws.Range("K129").Value = 2296.8
$ws.Range("M129").Value = 2703.2
$ws.Range("H131").Value = 12513861
$ws.Range("I131").Value = 17875748
$ws.Range("J131").Value = 2791.3333
$ws.Range("K131").Value = 53627244
$ws.Range("L131").Value = 8373.999899999999
$ws.Range("M131").Value = -53622204
$ws.Range("N131").Value = -18453.9999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 72000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 72000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 72000
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -72332
$ws.Range("H38").Value = 299146
$ws.Range("J38").Value = 299146
$ws.Range("L38").Value = 299146
$ws.Range("N38").Value = -300072
$ws.Range("H43").Value = 14329.25
$ws.Range("I43").Value = 4105.6665
$ws.Range("J43").Value = 45000
$ws.Range("K43").Value = 4105.6665
$ws.Range("L43").Value = 45000
$ws.Range("M43").Value = -3954.6665
$ws.Range("N43").Value = -45302
$ws.Range("H52").Value = 27999.5
$ws.Range("I52").Value = 15999
$ws.Range("J52").Value = 40000
$ws.Range("K52").Value = 15999
$ws.Range("L52").Value = 40000
$ws.Range("M52").Value = -15740
$ws.Range("N52").Value = -40518
$ws.Range("H132").Value = 2625.6
$ws.Range("I132").Value = 2491.75
$ws.Range("J132").Value = 2893.3
$ws.Range("K132").Value = 7475.25
$ws.Range("L132").Value = 8679.900000000001
$ws.Range("M132").Value = -4945.25
$ws.Range("N132").Value = -13739.9

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3550.3076
$ws.Range("J46").Value = 5081.5
$ws.Range("L46").Value = 5081.5
$ws.Range("N46").Value = -5457.5
$ws.Range("H61").Value = 1516.3334
$ws.Range("I61").Value = 1395
$ws.Range("K61").Value = 1395
$ws.Range("M61").Value = -1193
$ws.Range("H113").Value = 1516.3334
$ws.Range("I113").Value = 1395
$ws.Range("K113").Value = 1395
$ws.Range("M113").Value = 775
$ws.Range("H122").Value = 9273.5
$ws.Range("I122").Value = 7086.8184
$ws.Range("J122").Value = 10688.412
$ws.Range("K122").Value = 21260.4552
$ws.Range("L122").Value = 32065.236
$ws.Range("M122").Value = -18810.4552
$ws.Range("N122").Value = -36965.236
$ws.Range("H136").Value = 4751
$ws.Range("I136").Value = 1872.2
$ws.Range("K136").Value = 5616.6
$ws.Range("M136").Value = -3066.6

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872
$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360
$ws.Range("H113").Value = 1555.7059
$ws.Range("I113").Value = 1386.2727
$ws.Range("K113").Value = 4158.8181
$ws.Range("M113").Value = -1988.8181

